$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 5.685592333333333
$ws.Range("H2").Value = 17.056777
$ws.Range("I2").Value = 0.1188473284691575
$ws.Range("J2").Value = 0.1188473284691575
$ws.Range("M2").Value = 77.08952333333333
$ws.Range("N2").Value = 231.26857
$ws.Range("O2").Value = 0.2403816673726824
$ws.Range("P2").Value = 0.2403816673726824
$ws.Range("Q2").Value = 438.2996028443211
$ws.Range("R2").Value = 3944.69642559889
$ws.Range("S2").Value = 0.02856871898020496
$ws.Range("T2").Value = 0.02856871898020496
$ws.Range("G3").Value = 5.685592333333333
$ws.Range("H3").Value = 17.056777
$ws.Range("I3").Value = 0.1188473284691575
$ws.Range("J3").Value = 0.1188473284691575
$ws.Range("O3").Value = 0.3167483425780597
$ws.Range("P3").Value = 0.3167483425780597
$ws.Range("Q3").Value = 577.5426814821137
$ws.Range("R3").Value = 5197.884133339025
$ws.Range("S3").Value = 0.0376446943124359
$ws.Range("T3").Value = 0.0376446943124359
$ws.Range("G4").Value = 5.685592333333333
$ws.Range("H4").Value = 17.056777
$ws.Range("I4").Value = 0.1188473284691575
$ws.Range("J4").Value = 0.1188473284691575
$ws.Range("O4").Value = 0.4428699900492579
$ws.Range("P4").Value = 0.4428699900492579
$ws.Range("Q4").Value = 807.5064245615484
$ws.Range("R4").Value = 7267.557821053936
$ws.Range("S4").Value = 0.05263391517651669
$ws.Range("T4").Value = 0.05263391517651669
$ws.Range("I5").Value = 0.622926875404983
$ws.Range("J5").Value = 0.6229268754049829
$ws.Range("M5").Value = 77.08952333333333
$ws.Range("N5").Value = 231.26857
$ws.Range("O5").Value = 0.2403816673726824
$ws.Range("P5").Value = 0.2403816673726824
$ws.Range("Q5").Value = 2297.305337931197
$ws.Range("R5").Value = 20675.74804138078
$ws.Range("S5").Value = 0.149740200961105
$ws.Range("T5").Value = 0.149740200961105
$ws.Range("I6").Value = 0.622926875404983
$ws.Range("J6").Value = 0.6229268754049829
$ws.Range("O6").Value = 0.3167483425780597
$ws.Range("P6").Value = 0.3167483425780597
$ws.Range("S6").Value = 0.1973110553318579
$ws.Range("T6").Value = 0.1973110553318578
$ws.Range("I7").Value = 0.622926875404983
$ws.Range("J7").Value = 0.6229268754049829
$ws.Range("O7").Value = 0.4428699900492579
$ws.Range("P7").Value = 0.4428699900492579
$ws.Range("S7").Value = 0.2758756191120202
$ws.Range("T7").Value = 0.2758756191120201
$ws.Range("I8").Value = 0.2582257961258595
$ws.Range("J8").Value = 0.2582257961258594
$ws.Range("M8").Value = 77.08952333333333
$ws.Range("N8").Value = 231.26857
$ws.Range("O8").Value = 0.2403816673726824
$ws.Range("P8").Value = 0.2403816673726824
$ws.Range("Q8").Value = 952.3164327206111
$ws.Range("R8").Value = 8570.8478944855
$ws.Range("S8").Value = 0.06207274743137246
$ws.Range("T8").Value = 0.06207274743137244
$ws.Range("I9").Value = 0.2582257961258595
$ws.Range("J9").Value = 0.2582257961258594
$ws.Range("O9").Value = 0.3167483425780597
$ws.Range("P9").Value = 0.3167483425780597
$ws.Range("S9").Value = 0.08179259293376595
$ws.Range("T9").Value = 0.08179259293376592
$ws.Range("I10").Value = 0.2582257961258595
$ws.Range("J10").Value = 0.2582257961258594
$ws.Range("O10").Value = 0.4428699900492579
$ws.Range("P10").Value = 0.4428699900492579
$ws.Range("S10").Value = 0.1143604557607211
$ws.Range("T10").Value = 0.1143604557607211
